$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The survey rows for the 2025-08-10 batch (B1212:B1331) were actually
# collected on 2025-08-09 -- the date serial drops from 45879 to 45878.
$ws.Range("B1212:B1331").Value = 45878

# Restore the window/view state recorded for this edit: the sheet is
# zoomed out to 85% and the last-touched range (B1212:B1331, anchored at
# B1212) is selected.
$excel.ActiveWindow.Zoom = 85
$null = $ws.Range("B1212:B1331").Select()
